$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.133.93"
$ws.Range("E2").Value = "  -0.33%  "

$ws.Range("D3").Value = "3.460.16"
$ws.Range("E3").Value = "  -0.90%  "

$ws.Range("E4").Value = "  -0.01%  "

$ws.Range("D5").Value = "579.34"

$ws.Range("D6").Value = "149.36"
$ws.Range("E6").Value = "  +1.12%  "

$ws.Range("E7").Value = "  -0.04%  "

$ws.Range("D8").Value = "0.479"
$ws.Range("E8").Value = "  +0.11%  "

$ws.Range("E9").Value = "  +2.03%  "

$ws.Range("E10").Value = "  -2.02%  "

$ws.Range("D11").Value = "0.409"
$ws.Range("E11").Value = "  +2.20%  "

$ws.Range("D12").Value = "4.049.65"
$ws.Range("E12").Value = "  -0.96%  "

$ws.Range("E13").Value = "  +2.35%  "

$ws.Range("D14").Value = "28.57"
$ws.Range("E14").Value = "  -4.24%  "

$ws.Range("D15").Value = "3.455.34"
$ws.Range("E15").Value = "  -1.33%  "

$ws.Range("E16").Value = "  -1.38%  "

$ws.Range("D17").Value = "63.171.77"
$ws.Range("E17").Value = "  -0.36%  "

$ws.Range("D18").Value = "6.46"
$ws.Range("E18").Value = "  +2.84%  "

$ws.Range("D19").Value = "14.53"
$ws.Range("E19").Value = "  +1.31%  "

$ws.Range("D20").Value = "9.17"
$ws.Range("E20").Value = "  -3.36%  "

$ws.Range("D21").Value = "390.81"
$ws.Range("E21").Value = "  -0.68%  "

$ws.Range("D22").Value = "0.563"
$ws.Range("E22").Value = "  -0.66%  "

$ws.Range("D23").Value = "74.74"
$ws.Range("E23").Value = "  -1.03%  "

$ws.Range("E24").Value = "  +0.02%  "

$ws.Range("D25").Value = "3.591.95"
$ws.Range("E25").Value = "  -1.25%  "

$ws.Range("D26").Value = "0.0000115"
$ws.Range("E26").Value = "  -3.85%  "

$ws.Range("D27").Value = "0.183"
$ws.Range("E27").Value = "  -1.36%  "

$ws.Range("D28").Value = "7.67"
$ws.Range("E28").Value = "  -2.40%  "

$ws.Range("D29").Value = "1.00"
$ws.Range("E29").Value = "  +0.16%  "

$ws.Range("D30").Value = "8.08"
$ws.Range("E30").Value = "  -2.33%  "

$ws.Range("D31").Value = "2.12"
$ws.Range("E31").Value = "  -2.02%  "

$ws.Range("D32").Value = "0.999"

$ws.Range("B33").Value = "EthereumClassic"
$ws.Range("C33").Value = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D33").Value = "23.43"
$ws.Range("E33").Value = "  -1.84%  "

$ws.Range("B34").Value = "Fetch.AI"
$ws.Range("C34").Value = "https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet"
$ws.Range("D34").Value = "1.34"
$ws.Range("E34").Value = "  -5.68%  "

$ws.Range("B35").Value = "ImmutableX"
$ws.Range("C35").Value = "https://coinranking.com/coin/Z96jIvLU7+immutablex-imx"
$ws.Range("D35").Value = "1.63"
$ws.Range("E35").Value = "  +3.00%  "

$ws.Range("B36").Value = "NEARProtocol"
$ws.Range("C36").Value = "https://coinranking.com/coin/DCrsaMv68+nearprotocol-near"
$ws.Range("D36").Value = "5.36"
$ws.Range("E36").Value = "  -0.15%  "

$ws.Range("D37").Value = "32.00"
$ws.Range("E37").Value = "  -1.91%  "

$ws.Range("D38").Value = "7.05"
$ws.Range("E38").Value = "  -1.82%  "

$ws.Range("D39").Value = "170.08"
$ws.Range("E39").Value = "  -1.05%  "

$ws.Range("D40").Value = "3.496.54"
$ws.Range("E40").Value = "  -0.94%  "

$ws.Range("D41").Value = "0.0773"
$ws.Range("E41").Value = "  +0.44%  "

$ws.Range("E42").Value = "  -1.04%  "

$ws.Range("E43").Value = "  +0.69%  "

$ws.Range("E44").Value = "  -1.82%  "

$ws.Range("E45").Value = "  -1.66%  "

$ws.Range("E46").Value = "  -3.11%  "

$ws.Range("D47").Value = "2.583.70"
$ws.Range("E47").Value = "  -1.28%  "

$ws.Range("D48").Value = "2.29"
$ws.Range("E48").Value = "  -0.95%  "

$ws.Range("D49").Value = "6.91"
$ws.Range("E49").Value = "  +1.84%  "

$ws.Range("D50").Value = "22.65"
$ws.Range("E50").Value = "  -4.66%  "

$ws.Range("E51").Value = "  +0.00%  "
